$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long-form course names with their course codes.
$ws.Range("A2:A4").Value = "MT101"
$ws.Range("A5:A6").Value = "PH101"
$ws.Range("A7:A8").Value = "CS101"

# With the long descriptive text replaced by short course codes, the
# wrapped cells no longer need the extra height Excel had auto-fit
# earlier - restore the rows to the sheet's default height.
$ws.Rows("2:8").EntireRow.AutoFit()

# Move the active selection to A8.
$ws.Range("A8").Select()
